$p = $ppt.ActivePresentation

# .RGB setters in this host serialise red/blue swapped, so build the literal
# that round-trips to srgbClr val="274E13" (the deck's house green).
$greenColor = 0x134E27

# ---------------------------------------------------------------------------
# Slide 1: the green "Link to github Repository" / "Links to trello board"
# box. Append a hyperlinked "Github" run after the first line, and a
# hyperlinked "Trello Board" run after the third line.
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$body = $s1.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

# --- line 1: "Link to github Repository: " -> append hyperlinked "Github" --
$anchor1 = "Link to github Repository: "
$start1 = $tr.Text.IndexOf($anchor1) + 1
$line1 = $tr.Characters($start1, $anchor1.Length)
$line1.Text = $anchor1 + "Github"

$githubRun = $tr.Characters($start1 + $anchor1.Length, "Github".Length)
$githubRun.Font.Bold = $true
$githubRun.Font.Size = 20
$githubRun.Font.Color.RGB = $greenColor
$githubLink = $githubRun.ActionSettings(1).Hyperlink
$githubLink.Address = "https://github.com"

# --- line 3: "Links to trello board / project management tools:" -> add a --
# --- trailing space and append hyperlinked "Trello Board" -----------------
$anchor2 = "Links to trello board / project management tools:"
$start2 = $tr.Text.IndexOf($anchor2) + 1
$line3 = $tr.Characters($start2, $anchor2.Length)
$line3.Text = $anchor2 + " Trello Board"

$trelloRun = $tr.Characters($start2 + $anchor2.Length + 1, "Trello Board".Length)
$trelloRun.Font.Bold = $true
$trelloRun.Font.Size = 20
$trelloRun.Font.Color.RGB = $greenColor
$trelloLink = $trelloRun.ActionSettings(1).Hyperlink
$trelloLink.Address = "https://trello.com"

# ---------------------------------------------------------------------------
# Slide 15: title "Component 1 (Trello screenshot)" -- merge the two runs
# that make up the first line back into a single run (same visible text).
# ---------------------------------------------------------------------------
$s15 = $p.Slides.Item(15)
$title = $s15.Shapes.Item(1)
$ttr = $title.TextFrame.TextRange
$titleAnchor = "Component 1 (Trello screenshot)"
$titleStart = $ttr.Text.IndexOf($titleAnchor) + 1
$firstLine = $ttr.Characters($titleStart, $titleAnchor.Length)
$firstLine.Text = $titleAnchor
